# Updated cryptos list on Wed Apr 26 04:36:13 UTC 2023 with GitHub Actions
#
# Applies the per-cell text updates described by the target diff to the
# "cryptos" worksheet (Coin / Link / Price / Volume(1h) table).
#
# Notes:
#  - Columns D (Price) and E (Volume(1h)) hold values that are stored as
#    plain text in the workbook (t="inlineStr"/shared string), even though
#    many of them look like numbers (e.g. "339.00", "1.000", "0.4695").
#    Assigning such a string straight to .Value would make Excel silently
#    reinterpret it as a real number (dropping trailing zeros, losing the
#    original formatting). To keep these as text we prefix the risky ones
#    with a leading apostrophe, exactly like typing them in the Excel UI,
#    which forces Excel to store them as text (quote-prefixed) instead of
#    auto-converting them to numbers.
#  - Values that are safe as-is (URLs, coin names, multi-dot price strings
#    like "28.377.17", percentage strings with surrounding spaces like
#    "  +3.46%  ") are assigned directly since Excel already keeps them as
#    text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.377.17"
$ws.Range("E2").Value = "  +3.46%  "
$ws.Range("D3").Value = "1.868.75"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'339.00"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D7").Value = "'0.4695"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").Value = "'0.3971"
$ws.Range("E8").Value = "  +4.11%  "
$ws.Range("D9").Value = "'47.60"
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").Value = "'0.08022"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("E12").Value = "  +4.04%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.025"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.866.34"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "'7.253"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "'91.29"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "'0.06628"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "'17.56"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "28.390.94"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("D23").Value = "'5.475"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "2.104.44"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "'161.21"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "'19.76"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").Value = "'2.123"
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("E30").Value = "  +3.28%  "
$ws.Range("D31").Value = "'120.38"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("D32").Value = "'0.9709"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").Value = "'0.09501"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("D34").Value = "'3.596"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.375"
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'5.349"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "'0.06098"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").Value = "'0.02254"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("D39").Value = "'8.390"
$ws.Range("E39").Value = "  +3.85%  "
$ws.Range("D40").Value = "'1.185"
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("D41").Value = "'0.5951"
$ws.Range("E41").Value = "  +2.71%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'0.1875"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").Value = "'10.35"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("D45").Value = "'1.289"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("D46").Value = "'0.5583"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").Value = "'12.17"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").Value = "'1.960"
$ws.Range("E48").Value = "  +5.02%  "
$ws.Range("D49").Value = "'0.06870"
$ws.Range("E49").Value = "  +3.26%  "
$ws.Range("D50").Value = "'2.084"
$ws.Range("E50").Value = "  +18.23%  "
$ws.Range("D51").Value = "'111.72"
$ws.Range("E51").Value = "  +1.67%  "
